# Update price (D) and volume-change (E) columns for each coin row
# with the refreshed coinranking.com snapshot values. Column D holds
# free-form price text (not numbers), so force text format before
# assigning to avoid Excel auto-converting things like "0.9980" or
# "0.00001061" into a differently-formatted number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.890.02"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.769.91"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9980"
$ws.Range("E4").Value = "  -0.65%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.33"
$ws.Range("E5").Value = "  -1.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9986"
$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4268"
$ws.Range("E7").Value = "  -5.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3612"
$ws.Range("E8").Value = "  -3.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.58"
$ws.Range("E9").Value = "  -1.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07477"
$ws.Range("E10").Value = "  -4.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.108"
$ws.Range("E11").Value = "  -1.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9929"
$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.60"
$ws.Range("E13").Value = "  -1.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.130"
$ws.Range("E14").Value = "  -1.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.304"
$ws.Range("E15").Value = "  -1.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.789.85"
$ws.Range("E16").Value = "  +1.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.32"
$ws.Range("E17").Value = "  -0.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001061"
$ws.Range("E18").Value = "  -2.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06334"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9984"
$ws.Range("E20").Value = "  -0.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.22"
$ws.Range("E21").Value = "  -1.80%  "

$ws.Range("E22").Value = "  -4.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.851.75"
$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.37"
$ws.Range("E24").Value = "  -2.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.128"
$ws.Range("E25").Value = "  -9.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.49"
$ws.Range("E26").Value = "  +3.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.29"
$ws.Range("E27").Value = "  -2.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.991.58"
$ws.Range("E28").Value = "  +1.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.174"
$ws.Range("E29").Value = "  -7.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.77"
$ws.Range("E30").Value = "  -2.75%  "

$ws.Range("E31").Value = "  -4.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.690"
$ws.Range("E32").Value = "  -1.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08954"
$ws.Range("E33").Value = "  -3.77%  "

$ws.Range("E34").Value = "  -5.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.59"
$ws.Range("E35").Value = "  -1.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02319"
$ws.Range("E36").Value = "  -1.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6437"
$ws.Range("E37").Value = "  -1.14%  "

# Rows 38 and 39 swapped rank position (Algorand now ranks above
# InternetComputer); update Coin/Link/Price/Volume together.
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2111"
$ws.Range("E38").Value = "  -3.64%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.048"
$ws.Range("E39").Value = "  -1.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06050"
$ws.Range("E40").Value = "  -1.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.180"
$ws.Range("E41").Value = "  -1.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9967"
$ws.Range("E42").Value = "  -0.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.870"
$ws.Range("E43").Value = "  -2.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.393"
$ws.Range("E44").Value = "  -1.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.67"
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5968"
$ws.Range("E46").Value = "  -0.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.697"
$ws.Range("E47").Value = "  -1.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.15"
$ws.Range("E48").Value = "  -1.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.980"
$ws.Range("E49").Value = "  -1.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.144"
$ws.Range("E50").Value = "  -0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06889"
$ws.Range("E51").Value = "  -0.45%  "
